$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.490.89'
$ws.Range('E2').Value = '  +3.87%  '
$ws.Range('D3').Value = '2.276.63'
$ws.Range('E3').Value = '  +3.12%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '322.43'
$ws.Range('E5').Value = '  +2.29%  '
$ws.Range('D6').Value = '105.54'
$ws.Range('E6').Value = '  +6.85%  '
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.573'
$ws.Range('E9').Value = '  +2.48%  '
$ws.Range('E10').Value = '  +5.54%  '
$ws.Range('D11').Value = '0.0845'
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('D12').Value = '7.89'
$ws.Range('E12').Value = '  +3.00%  '
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').Value = '0.886'
$ws.Range('E14').Value = '  +3.47%  '
$ws.Range('D15').Value = '2.622.16'
$ws.Range('E15').Value = '  +3.28%  '
$ws.Range('D16').Value = '14.59'
$ws.Range('E16').Value = '  +2.76%  '
$ws.Range('D17').Value = '2.278.54'
$ws.Range('E17').Value = '  +3.32%  '
$ws.Range('D18').Value = '44.381.57'
$ws.Range('E18').Value = '  +3.86%  '
$ws.Range('D19').Value = '13.96'
$ws.Range('E19').Value = '  -2.95%  '
$ws.Range('D20').Value = '0.0000101'
$ws.Range('E20').Value = '  +4.93%  '
$ws.Range('D21').Value = '6.55'
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('D22').Value = '66.54'
$ws.Range('E22').Value = '  +2.16%  '
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('D24').Value = '240.47'
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('D25').Value = '2.24'
$ws.Range('E25').Value = '  +5.60%  '
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').Value = '10.24'
$ws.Range('D28').Value = '38.58'
$ws.Range('E28').Value = '  +12.89%  '
$ws.Range('D29').Value = '2.20'
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('E30').Value = '  +3.81%  '
$ws.Range('D31').Value = '20.73'
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0888'
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '162.34'
$ws.Range('E33').Value = '  +4.87%  '
$ws.Range('D34').Value = '2.78'
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('D35').Value = '0.118'
$ws.Range('E35').Value = '  +8.96%  '
$ws.Range('D36').Value = '2.03'
$ws.Range('E36').Value = '  +6.35%  '
$ws.Range('D37').Value = '3.15'
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('E38').Value = '  +0.71%  '
$ws.Range('D39').Value = '3.96'
$ws.Range('E39').Value = '  +5.14%  '
$ws.Range('D40').Value = '4.46'
$ws.Range('E40').Value = '  +1.18%  '
$ws.Range('D41').Value = '15.67'
$ws.Range('E41').Value = '  +28.01%  '
$ws.Range('D42').Value = '0.0330'
$ws.Range('E42').Value = '  +1.89%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').Value = '1.785.54'
$ws.Range('E44').Value = '  -3.51%  '
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('D46').Value = '86.60'
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').Value = '5.46'
$ws.Range('E47').Value = '  +2.33%  '
$ws.Range('D48').Value = '60.97'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').Value = '75.64'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = '1.73'
$ws.Range('E50').Value = '  +9.57%  '
$ws.Range('D51').Value = '104.57'
$ws.Range('E51').Value = '  +1.98%  '
